$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 362
$ws1.Range("F4").Value = 4761
$ws1.Range("F6").Value = 479

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 362
$ws4.Range("F4").Value = 4761
$ws4.Range("F8").Value = 479
